$d = $word.ActiveDocument

$d.Content.Find.Execute("84÷2=42, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2) | Out-Null
$d.Content.Find.Execute("12÷5=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2) | Out-Null
$d.Content.Find.Execute("84÷5=16, 4", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=13, 0", 2) | Out-Null
$d.Content.Find.Execute("80÷2=40, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷3=26, 1", 2) | Out-Null
$d.Content.Find.Execute("19÷7=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "52÷2=26, 0", 2) | Out-Null
$d.Content.Find.Execute("43÷8=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "91÷5=18, 1", 2) | Out-Null
$d.Content.Find.Execute("12÷7=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "51÷6=8, 3", 2) | Out-Null
$d.Content.Find.Execute("91÷8=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷7=10, 0", 2) | Out-Null
$d.Content.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "27÷3=9, 0", 2) | Out-Null
$d.Content.Find.Execute("81÷3=27, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷2=41, 0", 2) | Out-Null
$d.Content.Find.Execute("28÷7=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "20÷9=2, 2", 2) | Out-Null
$d.Content.Find.Execute("38÷5=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=10, 7", 2) | Out-Null
$d.Content.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "77÷9=8, 5", 2) | Out-Null
$d.Content.Find.Execute("83÷2=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "97÷8=12, 1", 2) | Out-Null
$d.Content.Find.Execute("80÷9=8, 8", $true, $false, $false, $false, $false, $true, 1, $false, "83÷5=16, 3", 2) | Out-Null
$d.Content.Find.Execute("81÷8=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "38÷9=4, 2", 2) | Out-Null
$d.Content.Find.Execute("51÷7=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "43÷6=7, 1", 2) | Out-Null
$d.Content.Find.Execute("26÷2=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷3=22, 2", 2) | Out-Null
$d.Content.Find.Execute("73÷5=14, 3", $true, $false, $false, $false, $false, $true, 1, $false, "77÷8=9, 5", 2) | Out-Null
$d.Content.Find.Execute("38÷7=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "14÷9=1, 5", 2) | Out-Null
$d.Content.Find.Execute("99÷8=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "17÷3=5, 2", 2) | Out-Null
$d.Content.Find.Execute("87÷4=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "14÷6=2, 2", 2) | Out-Null
$d.Content.Find.Execute("72÷3=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "99÷4=24, 3", 2) | Out-Null
$d.Content.Find.Execute("35÷7=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷7=7, 6", 2) | Out-Null
$d.Content.Find.Execute("99÷2=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "33÷2=16, 1", 2) | Out-Null
